$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Friday 2023-04-14 (row 5) attendance marking.
# Highlight the three attendee columns (Omar, Tommy, Faris) with the
# "Attended" green fill, and add a small legend in G5/H5 showing the
# two possible attendance states: green "Attended" and red "Absent".

$ws.Range("C5:E5").Interior.Color = 5287936   # RGB(0,176,80) green

$ws.Range("G5").Value = "Attended"
$ws.Range("G5").Interior.Color = 5287936      # RGB(0,176,80) green

$ws.Range("H5").Value = "Absent"
$ws.Range("H5").Interior.Color = 255          # RGB(255,0,0) red

[void]$ws.Range("G6").Select()
